$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 holds the all-digit string "11111111111". A plain .Value assignment
# would be auto-coerced to a number by Excel's smart-entry logic, and a
# leading apostrophe (quote-prefix) would stick a "stored as text" style on
# the cell. Neither matches the source file, so enter it as a text formula
# and immediately flatten it to a literal value (Copy / Paste Special
# Values), exactly like using Excel's Paste Values button - this keeps the
# cell as a genuine shared-string text value with the default style.
$ws.Range("A2").Formula = "=""11111111111"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("B2").Value = "Data dodania pacjenta do systemu:19.01.2023`n"
